# Update the NATMI ligand-receptor (Bmp15-Bmpr1b) sheet with refreshed TPM-based
# numbers: rows 2-3 (Sending cluster "FAPs" -> "ECs") get new edge-weight values,
# and four additional sender/target combinations (FAPs, MuSCs) are appended as
# rows 4-7, expanding the sheet from A1:T3 to A1:T7.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Bmp15"
$ws.Range("C2").Value = "Bmpr1b"
$ws.Range("D2").Value = "FAPs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.02770266666666667
$ws.Range("H2").Value = 0.083108
$ws.Range("I2").Value = 0.05180130905700151
$ws.Range("J2").Value = 0.05180130905700151
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 1.57938
$ws.Range("N2").Value = 4.73814
$ws.Range("O2").Value = 0.8416031693647025
$ws.Range("P2").Value = 0.8416031693647025
$ws.Range("Q2").Value = 0.04375303768
$ws.Range("R2").Value = 0.39377733912
$ws.Range("S2").Value = 0.04359614587961294
$ws.Range("T2").Value = 0.04359614587961294
# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Bmp15"
$ws.Range("C3").Value = "Bmpr1b"
$ws.Range("D3").Value = "MuSCs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.02770266666666667
$ws.Range("H3").Value = 0.083108
$ws.Range("I3").Value = 0.05180130905700151
$ws.Range("J3").Value = 0.05180130905700151
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.2972526666666667
$ws.Range("N3").Value = 0.891758
$ws.Range("O3").Value = 0.1583968306352975
$ws.Range("P3").Value = 0.1583968306352975
$ws.Range("Q3").Value = 0.008234691540444445
$ws.Range("R3").Value = 0.074112223864
$ws.Range("S3").Value = 0.00820516317738857
$ws.Range("T3").Value = 0.00820516317738857
# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Bmp15"
$ws.Range("C4").Value = "Bmpr1b"
$ws.Range("D4").Value = "FAPs"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.4126986666666667
$ws.Range("H4").Value = 1.238096
$ws.Range("I4").Value = 0.7717066171516261
$ws.Range("J4").Value = 0.7717066171516261
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 1.57938
$ws.Range("N4").Value = 4.73814
$ws.Range("O4").Value = 0.8416031693647025
$ws.Range("P4").Value = 0.8416031693647025
$ws.Range("Q4").Value = 0.65180802016
$ws.Range("R4").Value = 5.86627218144
$ws.Range("S4").Value = 0.6494707348145217
$ws.Range("T4").Value = 0.6494707348145217
# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Bmp15"
$ws.Range("C5").Value = "Bmpr1b"
$ws.Range("D5").Value = "MuSCs"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.4126986666666667
$ws.Range("H5").Value = 1.238096
$ws.Range("I5").Value = 0.7717066171516261
$ws.Range("J5").Value = 0.7717066171516261
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.2972526666666667
$ws.Range("N5").Value = 0.891758
$ws.Range("O5").Value = 0.1583968306352975
$ws.Range("P5").Value = 0.1583968306352975
$ws.Range("Q5").Value = 0.1226757791964445
$ws.Range("R5").Value = 1.104082012768
$ws.Range("S5").Value = 0.1222358823371045
$ws.Range("T5").Value = 0.1222358823371045
# Row 6
$ws.Range("A6").Value = "MuSCs"
$ws.Range("B6").Value = "Bmp15"
$ws.Range("C6").Value = "Bmpr1b"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 0.6666666666666666
$ws.Range("G6").Value = 0.09438566666666666
$ws.Range("H6").Value = 0.283157
$ws.Range("I6").Value = 0.1764920737913724
$ws.Range("J6").Value = 0.1764920737913724
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 1.57938
$ws.Range("N6").Value = 4.73814
$ws.Range("O6").Value = 0.8416031693647025
$ws.Range("P6").Value = 0.8416031693647025
$ws.Range("Q6").Value = 0.14907083422
$ws.Range("R6").Value = 1.34163750798
$ws.Range("S6").Value = 0.148536288670568
$ws.Range("T6").Value = 0.148536288670568
# Row 7
$ws.Range("A7").Value = "MuSCs"
$ws.Range("B7").Value = "Bmp15"
$ws.Range("C7").Value = "Bmpr1b"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 0.6666666666666666
$ws.Range("G7").Value = 0.09438566666666666
$ws.Range("H7").Value = 0.283157
$ws.Range("I7").Value = 0.1764920737913724
$ws.Range("J7").Value = 0.1764920737913724
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.2972526666666667
$ws.Range("N7").Value = 0.891758
$ws.Range("O7").Value = 0.1583968306352975
$ws.Range("P7").Value = 0.1583968306352975
$ws.Range("Q7").Value = 0.02805639111177777
$ws.Range("R7").Value = 0.252507520006
$ws.Range("S7").Value = 0.02795578512080444
$ws.Range("T7").Value = 0.02795578512080444
